$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new header row at the top with the text "review".
#    (this pushes every existing row down by one, which is what the
#    target sheet looks like: old A1..A18 become A2..A19)
# ------------------------------------------------------------------
$ws.Rows.Item(1).Insert() | Out-Null
$ws.Range("A1").Value = "review"

# ------------------------------------------------------------------
# 2. The row insert above does not move the worksheet's hyperlink
#    anchors, so re-point each existing hyperlink at its new
#    (shifted-down-by-one) cell. Hyperlinks.Add() forces Excel's
#    built-in "Hyperlink" cell style onto the target range, so we
#    snapshot + restore the original number/font formatting around
#    it and drop the now-unused "Hyperlink" style again.
# ------------------------------------------------------------------
$addrs = @("A4", "A5", "A6", "A15", "A16")
$urls  = @(
    "https://t.co/9MbNWxgvYC?amp=1",
    "https://ift.tt/355IO28",
    "https://drive.google.com/file/d/1MRRL70cq3D-GqqF2cKSw22iatmWKtl9n/view?usp=sharing",
    "http://bit.ly/2MSnkiY",
    "https://punto-informatico.it/adobe-cc-libraries-microsoft-office-365/"
)

for ($i = 0; $i -lt $addrs.Length; $i++) {
    $ws.Range($addrs[$i]).Copy() | Out-Null
    $ws.Cells.Item(1, 26 + $i).PasteSpecial(-4122) | Out-Null
}

# Drop the stale (pre-shift) hyperlink anchors before re-adding them at
# their new locations.
$ws.Hyperlinks.Delete() | Out-Null

for ($i = 0; $i -lt $addrs.Length; $i++) {
    $ws.Hyperlinks.Add($ws.Range($addrs[$i]), $urls[$i]) | Out-Null
}

$wb.Styles.Item("Hyperlink").Delete() | Out-Null

for ($i = 0; $i -lt $addrs.Length; $i++) {
    $ws.Cells.Item(1, 26 + $i).Copy() | Out-Null
    $ws.Range($addrs[$i]).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item(1, 26 + $i).Clear() | Out-Null
}

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3. Restore the selection to where the author left it.
# ------------------------------------------------------------------
$ws.Range("A11").Select() | Out-Null
